# Seating-chart standby drag fix
#
# Olivia Gonzalez (Contestants row 17) moves down to the bottom of the
# "available" block (new row 32); every row in between (old rows 18-32)
# shifts up by one. Her matching standby record (Standbys row 3) is
# removed entirely.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Contestants sheet: relocate row 17 -> row 32
# ---------------------------------------------------------------------
$contestants = $wb.Worksheets.Item("Contestants")

# Stash the row-17 record (columns A:J only - K:M are blank for this
# record) in an unused scratch row so it survives the reshuffle below.
$contestants.Range("A17:J17").Copy($contestants.Range("A400:J400"))
$excel.CutCopyMode = $false

# Removing row 17 shifts rows 18-32 up into 17-31.
$contestants.Rows("17").Delete()

# Make room for the record at its new home (row 32, right after the
# last "available" contestant and before the next block).
$contestants.Rows("32").Insert()

# Drop the stashed record into the freshly inserted row, then clear the
# scratch copy.
$contestants.Range("A400:J400").Copy($contestants.Range("A32:J32"))
$excel.CutCopyMode = $false
$contestants.Range("A400:J400").Clear()

# ---------------------------------------------------------------------
# 2) Standbys sheet: drop the now-obsolete standby row for her
# ---------------------------------------------------------------------
$standbys = $wb.Worksheets.Item("Standbys")
$standbys.Rows("3").Delete()
